# Weekly fruit/vegetable price update: insert two new daily-price rows
# (row 240 and 241) into the "Zapallo italiano" sheet, pushing the
# previously-existing rows 240-296 down to 242-298.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 240 (shifts rows 240:296 -> 242:298)
$ws.Range("A240:R241").EntireRow.Insert()

# New row 240 data
$ws.Cells.Item(240, 1).Value = 1
$ws.Cells.Item(240, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(240, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(240, 4).Value = 44642
$ws.Cells.Item(240, 5).Value = 15
$ws.Cells.Item(240, 6).Value = 100112032
$ws.Cells.Item(240, 7).Value = "Zapallo italiano"
$ws.Cells.Item(240, 8).Value = "Huracán"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 120
$ws.Cells.Item(240, 11).Value = 7000
$ws.Cells.Item(240, 12).Value = 7500
$ws.Cells.Item(240, 13).Value = 7250
$ws.Cells.Item(240, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(240, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(240, 16).Value = 104
$ws.Cells.Item(240, 17).Value = 70
$ws.Cells.Item(240, 18).Value = "Hortaliza"

# New row 241 data
$ws.Cells.Item(241, 1).Value = 1
$ws.Cells.Item(241, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(241, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(241, 4).Value = 44642
$ws.Cells.Item(241, 5).Value = 15
$ws.Cells.Item(241, 6).Value = 100112032
$ws.Cells.Item(241, 7).Value = "Zapallo italiano"
$ws.Cells.Item(241, 8).Value = "Huracán"
$ws.Cells.Item(241, 9).Value = "Segunda"
$ws.Cells.Item(241, 10).Value = 140
$ws.Cells.Item(241, 11).Value = 5500
$ws.Cells.Item(241, 12).Value = 6000
$ws.Cells.Item(241, 13).Value = 5750
$ws.Cells.Item(241, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(241, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(241, 16).Value = 58
$ws.Cells.Item(241, 17).Value = 100
$ws.Cells.Item(241, 18).Value = "Hortaliza"
